# "Actualización 10 de Mayo"
# Updates statistics figures on the three parcial sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "1er Parcial" ----
$ws = $wb.Worksheets.Item("1er Parcial")

# Row 3
$ws.Range("I3").Value = 6.4
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

# Row 4
$ws.Range("I4").Value = 6.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0

# Row 5
$ws.Range("I5").Value = 6.2
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3.45

# Row 12
$ws.Range("E12").Value = 17
$ws.Range("F12").Value = 11
$ws.Range("G12").Value = 60.71
$ws.Range("H12").Value = 39.29
$ws.Range("I12").Value = 7.2
$ws.Range("J12").Value = 11
$ws.Range("K12").Value = 39.29

# ---- Sheet "2o Parcial" ----
$ws = $wb.Worksheets.Item("2o Parcial")

# Row 2
$ws.Range("E2").Value = 28
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 77.78
$ws.Range("H2").Value = 22.22
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 16.67

# Row 5
$ws.Range("E5").Value = 18
$ws.Range("F5").Value = 11
$ws.Range("G5").Value = 62.07
$ws.Range("H5").Value = 37.93
$ws.Range("I5").Value = 6.7
$ws.Range("J5").Value = 8
$ws.Range("K5").Value = 27.59

# ---- Sheet "3er Parcial" ----
$ws = $wb.Worksheets.Item("3er Parcial")

# Row 2
$ws.Range("E2").Value = 28
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 77.78
$ws.Range("H2").Value = 22.22
$ws.Range("I2").Value = 6.5

# Row 3
$ws.Range("I3").Value = 6.1
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

# Row 4
$ws.Range("I4").Value = 6.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0

# Row 5
$ws.Range("I5").Value = 6.4
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3.45

# Row 12
$ws.Range("E12").Value = 17
$ws.Range("F12").Value = 11
$ws.Range("G12").Value = 60.71
$ws.Range("H12").Value = 39.29
$ws.Range("I12").Value = 7.2
$ws.Range("J12").Value = 11
$ws.Range("K12").Value = 39.29

Write-Host "Updated 1er/2o/3er Parcial sheets (Actualización 10 de Mayo)"
